$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F49").Value = 3
$ws.Range("G49").Value = 2750.82
$ws.Range("B52").Value = 3536.05
$ws.Range("F64").Value = 110
$ws.Range("G64").Value = 8928.700000000001
$ws.Range("F73").Value = 66
$ws.Range("G73").Value = 5212.68
$ws.Range("F74").Value = 140
$ws.Range("G74").Value = 19649
$ws.Range("F77").Value = 233
$ws.Range("G77").Value = 10890.42
$ws.Range("F79").Value = 71
$ws.Range("G79").Value = 4414.07
$ws.Range("F81").Value = 9
$ws.Range("G81").Value = 275.22
$ws.Range("B90").Value = 165146.23
$ws.Range("F115").Value = 176
$ws.Range("G115").Value = 17038.56
$ws.Range("B117").Value = 10879.62
$ws.Range("F149").Value = 211
$ws.Range("G149").Value = 13672.8
$ws.Range("B156").Value = 28473.03
$ws.Range("B219").Value = 61610
$ws.Range("E219").Value = 122.71
$ws.Range("F219").Value = -58
$ws.Range("G219").Value = -5957.18
$ws.Range("B220").Value = 63565
$ws.Range("E220").Value = 109.19
$ws.Range("F220").Value = 60
$ws.Range("G220").Value = 6162.6
$ws.Range("B232").Value = 63510
$ws.Range("E232").Value = 50.66
$ws.Range("F232").Value = 113
$ws.Range("G232").Value = 5383.32
$ws.Range("B233").Value = 55356
$ws.Range("E233").Value = 54.04
$ws.Range("F233").Value = -158
$ws.Range("G233").Value = -7527.12
$ws.Range("B243").Value = 63560
$ws.Range("E243").Value = 134.87
$ws.Range("F243").Value = 1
$ws.Range("G243").Value = 126.86
$ws.Range("B244").Value = 60325
$ws.Range("E244").Value = 151.57
$ws.Range("F244").Value = -102
$ws.Range("G244").Value = -12939.72
$ws.Range("F255").Value = 508
$ws.Range("G255").Value = 87035.64
$ws.Range("B260").Value = 165441.98
$ws.Range("F292").Value = 37
$ws.Range("G292").Value = 3080.99
$ws.Range("F293").Value = 25
$ws.Range("G293").Value = 1758
$ws.Range("F294").Value = 22
$ws.Range("G294").Value = 1569.92
$ws.Range("B304").Value = 160961.89
$ws.Range("B322").Value = 58047
$ws.Range("D322").Value = 105.54
$ws.Range("E322").Value = 126.1
$ws.Range("F322").Value = 39
$ws.Range("G322").Value = 4116.06
$ws.Range("B323").Value = 47097
$ws.Range("D323").Value = 112.28
$ws.Range("E323").Value = 134.16
$ws.Range("F323").Value = 15
$ws.Range("G323").Value = 1684.2
$ws.Range("F326").Value = 57
$ws.Range("G326").Value = 1695.18
$ws.Range("B330").Value = 25128.55
$ws.Range("F334").Value = 186
$ws.Range("G334").Value = 9638.52
$ws.Range("B346").Value = 23286.48
$ws.Range("B364").Value = 65068
$ws.Range("E364").Value = 13.97
$ws.Range("F364").Value = 63
$ws.Range("G364").Value = 828.45
$ws.Range("B365").Value = 53602
$ws.Range("E365").Value = 15.69
$ws.Range("F365").Value = -231
$ws.Range("G365").Value = -3037.65
$ws.Range("B366").Value = 53263
$ws.Range("E366").Value = 15.29
$ws.Range("F366").Value = -309
$ws.Range("G366").Value = -3958.29
$ws.Range("B367").Value = 65066
$ws.Range("E367").Value = 13.61
$ws.Range("F367").Value = 90
$ws.Range("G367").Value = 1152.9
$ws.Range("B372").Value = 64922
$ws.Range("E372").Value = 20.98
$ws.Range("F372").Value = 67
$ws.Range("G372").Value = 1321.91
$ws.Range("B373").Value = 45706
$ws.Range("E373").Value = 23.58
$ws.Range("F373").Value = -202
$ws.Range("G373").Value = -3985.46
$ws.Range("B375").Value = 64927
$ws.Range("E375").Value = 17.26
$ws.Range("F375").Value = 106
$ws.Range("G375").Value = 1719.32
$ws.Range("B376").Value = 45718
$ws.Range("E376").Value = 19.38
$ws.Range("F376").Value = -294
$ws.Range("G376").Value = -4768.68
$ws.Range("B380").Value = 45709
$ws.Range("E380").Value = 15.69
$ws.Range("F380").Value = -300
$ws.Range("G380").Value = -3945
$ws.Range("B381").Value = 64925
$ws.Range("E381").Value = 13.97
$ws.Range("F381").Value = 111
$ws.Range("G381").Value = 1459.65
$ws.Range("B382").Value = 64919
$ws.Range("E382").Value = 27.97
$ws.Range("F382").Value = 61
$ws.Range("G382").Value = 1604.3
$ws.Range("B383").Value = 45702
$ws.Range("E383").Value = 31.43
$ws.Range("F383").Value = -215
$ws.Range("G383").Value = -5654.5
$ws.Range("B385").Value = 53595
$ws.Range("E385").Value = 17.61
$ws.Range("F385").Value = -335
$ws.Range("G385").Value = -4934.55
$ws.Range("B386").Value = 65067
$ws.Range("E386").Value = 15.65
$ws.Range("F386").Value = 126
$ws.Range("G386").Value = 1855.98
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79
$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 105
$ws.Range("G474").Value = 3447.15
$ws.Range("B572").Value = 65362
$ws.Range("F572").Value = 18
$ws.Range("G572").Value = 735.66
$ws.Range("B573").Value = 65079
$ws.Range("F573").Value = 6
$ws.Range("G573").Value = 245.22
$ws.Range("F599").Value = 1280
$ws.Range("G599").Value = 208780.8
$ws.Range("F602").Value = 305
$ws.Range("G602").Value = 44118.25
$ws.Range("B606").Value = 354165.95
$ws.Range("B619").Value = 1549427.95
$ws.Range("B620").Value = 1549427.95
